# ---------------------------------------------------------------------------
# etl/source/metadata.xlsx -- "includes Demography Interpolated file AND
# variants are not seperated by underscore char"
#
# 1. Adds a new "Include" column (F) with a header cell + a literal 1 in
#    every existing data row (2-80).
# 2. Appends a new data row (81) for
#    WPP2015_INT_F01_ANNUAL_DEMOGRAPHIC_INDICATORS.XLS / NewFormat, carrying
#    the new "Include" flag too.
# 3. Updates the current selection to E75 (the last edited cell).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New "Include" column -----------------------------------------------
$ws.Range("F1").Value = "Include"

for ($r = 2; $r -le 80; $r++) {
    $ws.Cells.Item($r, 6).Value = 1
}

# Match the header's look (font/fill) to the rest of row 1 -- copy the
# existing header formatting onto the new F1 cell (data cells F2:F80 keep
# the workbook's default/no style, same as the rest of the table body).
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- 2. New row for the Demography Interpolated file -----------------------
# Fill the new row's cells in the same order the source workbook's shared
# string table was built in, so the underlying values line up exactly.
$ws.Range("B81").Value = "NewFormat"
$ws.Range("D81").Value = "NewFormat"
$ws.Range("E81").Value = "ddf--datapoints--NewFormat--by--ref_area_code--year--variant-"
$ws.Range("A81").Value = "WPP2015_INT_F01_ANNUAL_DEMOGRAPHIC_INDICATORS.XLS"
$ws.Range("C81").Value = "na"
$ws.Range("F81").Value = 1

# A81:D81 pick up the same styled look as the rest of the table (A-D cols);
# E81/F81 stay on the default/no style, matching the existing row pattern.
$ws.Range("A1:D1").Copy()
$ws.Range("A81:D81").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- 3. View bookkeeping -----------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 47
$win.ScrollColumn = 1
$win.Left = 0
$win.Top = 460
$win.Width = 38400
$win.Height = 19540

$ws.Range("E75").Select()
